# Updated cryptos list on Fri Nov  3 10:16:20 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $value) {
    # Prefix with an apostrophe so Excel stores the value as text even when
    # it looks numeric, then reset the style so no extra number-format /
    # quote-prefix style gets attached to the cell.
    $ws.Range($cellAddr).Value = "'" + $value
    $ws.Range($cellAddr).Style = "Normal"
}

# Rows 13-15 changed coin identity/order (B, C, D, E all change)
Set-TextValue "B13" "WrappedEther"
Set-TextValue "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.835.87"
Set-TextValue "E13" "  -0.19%  "

Set-TextValue "B14" "Chainlink"
Set-TextValue "C14" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D14" "11.07"
Set-TextValue "E14" "  -1.66%  "

Set-TextValue "B15" "Polygon"
Set-TextValue "C15" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D15" "0.657"
Set-TextValue "E15" "  -2.74%  "

# Remaining rows: only Price (D) and/or Volume(1h) (E) updated
Set-TextValue "D2" "34.494.08"
Set-TextValue "E2" "  -2.91%  "

Set-TextValue "D3" "1.801.26"
Set-TextValue "E3" "  -2.17%  "

Set-TextValue "D5" "229.37"
Set-TextValue "E5" "  -1.10%  "

Set-TextValue "D6" "0.610"
Set-TextValue "E6" "  -1.45%  "

Set-TextValue "E7" "  +0.62%  "

Set-TextValue "D8" "38.83"
Set-TextValue "E8" "  -10.92%  "

Set-TextValue "D9" "0.319"
Set-TextValue "E9" "  +2.16%  "

Set-TextValue "D10" "0.0675"
Set-TextValue "E10" "  -3.88%  "

Set-TextValue "D11" "0.0988"
Set-TextValue "E11" "  -2.10%  "

Set-TextValue "D12" "2.065.07"
Set-TextValue "E12" "  -2.00%  "

Set-TextValue "D16" "4.54"
Set-TextValue "E16" "  -4.15%  "

Set-TextValue "D17" "34.534.37"
Set-TextValue "E17" "  -2.70%  "

Set-TextValue "D18" "68.85"
Set-TextValue "E18" "  -2.30%  "

Set-TextValue "D19" "243.31"
Set-TextValue "E19" "  -0.68%  "

Set-TextValue "D20" "0.0₃0776"
Set-TextValue "E20" "  -3.18%  "

Set-TextValue "D21" "11.71"
Set-TextValue "E21" "  -2.84%  "

Set-TextValue "D22" "4.65"
Set-TextValue "E22" "  -2.18%  "

Set-TextValue "E23" "  +0.51%  "

Set-TextValue "D24" "2.22"
Set-TextValue "E24" "  -0.04%  "

Set-TextValue "D25" "171.86"
Set-TextValue "E25" "  +0.10%  "

Set-TextValue "D26" "7.68"
Set-TextValue "E26" "  -3.68%  "

Set-TextValue "D27" "17.08"
Set-TextValue "E27" "  -4.30%  "

Set-TextValue "E28" "  -1.39%  "

Set-TextValue "D29" "1.47"
Set-TextValue "E29" "  -7.25%  "

Set-TextValue "E30" "  +0.52%  "

Set-TextValue "D31" "4.02"
Set-TextValue "E31" "  +2.03%  "

Set-TextValue "D32" "0.0538"
Set-TextValue "E32" "  -2.56%  "

Set-TextValue "D33" "3.86"
Set-TextValue "E33" "  -5.49%  "

Set-TextValue "E34" "  +9.56%  "

Set-TextValue "D35" "1.77"
Set-TextValue "E35" "  -4.22%  "

Set-TextValue "D36" "0.688"
Set-TextValue "E36" "  -0.64%  "

Set-TextValue "D37" "90.59"
Set-TextValue "E37" "  -5.27%  "

Set-TextValue "E38" "  +4.28%  "

Set-TextValue "D39" "1.316.51"
Set-TextValue "E39" "  -2.60%  "

Set-TextValue "D40" "0.0190"
Set-TextValue "E40" "  -2.85%  "

Set-TextValue "D41" "2.45"
Set-TextValue "E41" "  -0.19%  "

Set-TextValue "D42" "0.950"
Set-TextValue "E42" "  -6.58%  "

Set-TextValue "D43" "14.20"
Set-TextValue "E43" "  -8.36%  "

Set-TextValue "D44" "2.70"
Set-TextValue "E44" "  -4.12%  "

Set-TextValue "D45" "2.18"
Set-TextValue "E45" "  -11.69%  "

Set-TextValue "D46" "6.16"
Set-TextValue "E46" "  -2.16%  "

Set-TextValue "E47" "  -0.84%  "

Set-TextValue "D48" "1.986.34"
Set-TextValue "E48" "  -1.22%  "

Set-TextValue "E49" "  +0.56%  "

Set-TextValue "D50" "0.0660"
Set-TextValue "E50" "  +4.10%  "

Set-TextValue "D51" "97.17"
Set-TextValue "E51" "  -5.81%  "
